$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values per repulled data
$ws.Range("F4").Value = 7
$ws.Range("F8").Value = 1
$ws.Range("F12").Value = 8
$ws.Range("F13").Value = 9
$ws.Range("F17").Value = -3
$ws.Range("F18").Value = -1
